$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'242.94"
$ws.Range("D3").Value = "'23.04"
$ws.Range("D4").Value = "'5.409"
$ws.Range("D5").Value = "'0.05898"
$ws.Range("D6").Value = "'3.439"
$ws.Range("D7").Value = "'6.523"
$ws.Range("D9").Value = "'0.9382"
$ws.Range("D10").Value = "'0.1423"
$ws.Range("D11").Value = "'0.07424"
$ws.Range("D12").Value = "'0.03294"
$ws.Range("D13").Value = "'0.03065"
$ws.Range("D15").Value = "'3.856"
$ws.Range("D16").Value = "'0.001581"
$ws.Range("D17").Value = "'0.04679"
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D18").Value = "'0.006008"
$ws.Range("E18").Value = "17TigerCashTCH"
$ws.Range("B19").Value = "BitKan"
$ws.Range("C19").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D19").Value = "'0.001258"
$ws.Range("E19").Value = "18BitKanKAN"
$ws.Range("B20").Value = "HotbitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D20").Value = "'0.004906"
$ws.Range("E20").Value = "19HotbitTokenHTB"
$ws.Range("B21").Value = "NitroEx"
$ws.Range("C21").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D21").Value = "'0.00006801"
$ws.Range("E21").Value = "20NitroExNTX"
$ws.Range("B22").Value = "LEO"
$ws.Range("C22").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D22").Value = "'3.562"
$ws.Range("E22").Value = "21LEOLEO"
$ws.Range("B23").Value = "BTSEToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D23").Value = "'2.126"
$ws.Range("E23").Value = "22BTSETokenBTSE"
$ws.Range("B24").Value = "One"
$ws.Range("C24").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D24").Value = "'0.01120"
$ws.Range("E24").Value = "23OneONEBestin24h"
$ws.Range("D27").Value = "'0.0003700"
$ws.Range("D40").Value = "'0.03966"
$ws.Range("D41").Value = "'0.006188"
$ws.Range("D42").Value = "'0.1073"
$ws.Range("D43").Value = "'0.003001"
$ws.Range("D44").Value = "'0.009499"
$ws.Range("E44").Value = "43LocalTradersLCT"
$ws.Range("D45").Value = "'0.00005227"
$ws.Range("D47").Value = "'0.6702"
$ws.Range("D48").Value = "'0.002327"
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("D50").Value = "'0.0002000"
